# Updates cryptos list data (Coin/Link/Price/Volume columns) on sheet1,
# matching the "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.228.54"
$ws.Range("E2").Value = "  -0.96%  "

# Row 3
$ws.Range("D3").Value = "1.860.06"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9986"
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.19"
$ws.Range("E5").Value = "  -2.58%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9986"
$ws.Range("E6").Value = "  -0.02%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4700"
$ws.Range("E7").Value = "  -1.73%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2816"
$ws.Range("E8").Value = "  -0.95%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06555"
$ws.Range("E9").Value = "  -2.40%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.12"
$ws.Range("E10").Value = "  +3.47%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07797"
$ws.Range("E11").Value = "  +0.36%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.34"
$ws.Range("E12").Value = "  -6.36%  "

# Row 13
$ws.Range("D13").Value = "1.860.04"
$ws.Range("E13").Value = "  -2.91%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.111"
$ws.Range("E14").Value = "  -2.24%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6664"
$ws.Range("E15").Value = "  -1.22%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "284.82"
$ws.Range("E16").Value = "  -3.18%  "

# Row 17
$ws.Range("D17").Value = "30.248.08"
$ws.Range("E17").Value = "  -0.93%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9996"
$ws.Range("E18").Value = "  +0.07%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.452"
$ws.Range("E19").Value = "  +0.77%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.59"
$ws.Range("E20").Value = "  -0.75%  "

# Row 21
$ws.Range("D21").Value = "2.105.98"
$ws.Range("E21").Value = "  -2.54%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007247"
$ws.Range("E22").Value = "  -3.34%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9981"
$ws.Range("E23").Value = "  +0.15%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.149"
$ws.Range("E24").Value = "  -2.63%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.53"
$ws.Range("E25").Value = "  +0.19%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.314"
$ws.Range("E26").Value = "  -0.95%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.06"
$ws.Range("E27").Value = "  -3.13%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.919"
$ws.Range("E28").Value = "  -8.47%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.340"
$ws.Range("E29").Value = "  -3.29%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09599"
$ws.Range("E30").Value = "  -3.57%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.422"
$ws.Range("E31").Value = "  -3.60%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.469"
$ws.Range("E32").Value = "  -3.11%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.105"
$ws.Range("E33").Value = "  -3.77%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04677"
$ws.Range("E34").Value = "  -1.37%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.102"
$ws.Range("E35").Value = "  -1.20%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6999"
$ws.Range("E36").Value = "  -3.95%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9986"
$ws.Range("E37").Value = "  +0.16%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.709"
$ws.Range("E38").Value = "  +0.10%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01855"
$ws.Range("E39").Value = "  -3.11%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.378"
$ws.Range("E40").Value = "  +0.35%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.509"
$ws.Range("E41").Value = "  -4.16%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.17"
$ws.Range("E42").Value = "  -3.45%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8562"
$ws.Range("E43").Value = "  -1.31%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.940"
$ws.Range("E44").Value = "  -1.12%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.09"
$ws.Range("E45").Value = "  -1.98%  "

# Row 46
$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4168"
$ws.Range("E46").Value = "  -2.47%  "

# Row 47
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9987"
$ws.Range("E47").Value = "  +0.05%  "

# Row 48
$ws.Range("D48").Value = "1.010.07"
$ws.Range("E48").Value = "  +4.69%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.211"
$ws.Range("E49").Value = "  -3.06%  "

# Row 50
$ws.Range("E50").Value = "  +2.69%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.75"
$ws.Range("E51").Value = "  -2.73%  "
